$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E holds text-typed (inline-string) estimate values that merely look
# numeric/percent. Mark the range as Text before writing so the engine keeps
# them as strings (not auto-coerced numbers), then clear the format again so
# no stray number-format style sticks to the cell (matches original styling).
$eCells = @("E2","E3","E4","E5","E9","E10","E11","E12","E13","E14","E15","E16","E17","E18","E19","E20","E21","E22","E23","E24","E25","E28")
foreach ($addr in $eCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("E2").Value = "689.4"
$ws.Range("F2").Value = 634.6
$ws.Range("G2").Value = 743.1

# Row 3
$ws.Range("E3").Value = "44.3"
$ws.Range("F3").Value = -24.4
$ws.Range("G3").Value = 113

# Row 4
$ws.Range("E4").Value = "-27.4"
$ws.Range("F4").Value = -92.90000000000001
$ws.Range("G4").Value = 44.2

# Row 5
$ws.Range("E5").Value = "-19.1"
$ws.Range("F5").Value = -118.7
$ws.Range("G5").Value = 81.40000000000001

# Row 9
$ws.Range("E9").Value = "596.5"
$ws.Range("F9").Value = 504.3
$ws.Range("G9").Value = 686.6

# Row 10
$ws.Range("E10").Value = "250.1"
$ws.Range("F10").Value = 101
$ws.Range("G10").Value = 406.2

# Row 11
$ws.Range("E11").Value = "85"
$ws.Range("F11").Value = -44
$ws.Range("G11").Value = 215.9

# Row 12
$ws.Range("E12").Value = "-161.4"
$ws.Range("F12").Value = -369.5
$ws.Range("G12").Value = 43.6

# Row 13
$ws.Range("E13").Value = "19%"

# Row 14
$ws.Range("E14").Value = "4%"

# Row 15
$ws.Range("E15").Value = "77%"

# Row 16
$ws.Range("E16").Value = "702.6"
$ws.Range("F16").Value = 649.3
$ws.Range("G16").Value = 758.2

# Row 17
$ws.Range("E17").Value = "-19.7"
$ws.Range("F17").Value = -80.3
$ws.Range("G17").Value = 37.9

# Row 18
$ws.Range("E18").Value = "-12.5"
$ws.Range("F18").Value = -54.8
$ws.Range("G18").Value = 30.3

# Row 19
$ws.Range("E19").Value = "27.7"
$ws.Range("G19").Value = 66.59999999999999

# Row 20
$ws.Range("E20").Value = "17.8"
$ws.Range("F20").Value = -19.5
$ws.Range("G20").Value = 55.6

# Row 21
$ws.Range("E21").Value = "42.9"
$ws.Range("F21").Value = -11.8
$ws.Range("G21").Value = 96.2

# Row 22
$ws.Range("E22").Value = "-2.6"
$ws.Range("F22").Value = -62
$ws.Range("G22").Value = 56.9

# Row 23
$ws.Range("E23").Value = "-1.2"
$ws.Range("F23").Value = -55.7
$ws.Range("G23").Value = 54.4

# Row 24
$ws.Range("E24").Value = "33.7"
$ws.Range("F24").Value = -22.6
$ws.Range("G24").Value = 89.3

# Row 25
$ws.Range("E25").Value = "-78.6"
$ws.Range("F25").Value = -160.4
$ws.Range("G25").Value = 5.4

# Row 28
$ws.Range("E28").Value = "31%"

foreach ($addr in $eCells) {
    $ws.Range($addr).ClearFormats()
}
